$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.409.59"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.675.82"
$ws.Range("E3").Value = "  -0.63%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "642.01"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.62"
$ws.Range("E6").Value = "  -0.90%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  +0.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.144"
$ws.Range("E9").Value = "  -1.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.06"
$ws.Range("E10").Value = "  -1.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.444"
$ws.Range("E11").Value = "  +0.56%  "

$ws.Range("E12").Value = "  -1.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.296.25"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.46"
$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.679.72"
$ws.Range("E15").Value = "  -0.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.381.77"
$ws.Range("E16").Value = "  -0.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.117"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.91"
$ws.Range("E18").Value = "  -1.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.44"
$ws.Range("E19").Value = "  -0.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "466.14"
$ws.Range("E20").Value = "  -0.98%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.84"
$ws.Range("E21").Value = "  -1.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.645"
$ws.Range("E22").Value = "  -1.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.25"
$ws.Range("E23").Value = "  -1.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.822.40"
$ws.Range("E24").Value = "  -0.51%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000124"
$ws.Range("E26").Value = "  -0.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.76"
$ws.Range("E27").Value = "  -2.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.95"
$ws.Range("E28").Value = "  -2.36%  "

$ws.Range("E29").Value = "  -3.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.69"
$ws.Range("E30").Value = "  -2.71%  "

$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.99"
$ws.Range("E32").Value = "  -1.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.79"
$ws.Range("E33").Value = "  -0.96%  "

$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.164"
$ws.Range("E34").Value = "  +2.78%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.40"
$ws.Range("E35").Value = "  -2.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.668.83"
$ws.Range("E36").Value = "  -0.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.39"
$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.86"
$ws.Range("E39").Value = "  -6.50%  "

$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.22"
$ws.Range("E41").Value = "  -2.72%  "

$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "176.59"
$ws.Range("E42").Value = "  +4.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0892"
$ws.Range("E43").Value = "  -1.87%  "

$ws.Range("E44").Value = "  -2.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.78"
$ws.Range("E45").Value = "  -0.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.69"
$ws.Range("E46").Value = "  -1.87%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.23"
$ws.Range("E47").Value = "  -5.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.06"
$ws.Range("E48").Value = "  -4.16%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.80"
$ws.Range("E49").Value = "  -0.94%  "

$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.24"
$ws.Range("E50").Value = "  -4.68%  "

$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000264"
$ws.Range("E51").Value = "  -6.38%  "
